# Updates cryptos list prices/volumes and reorders two coin pairs,
# matching the "Updated cryptos list ... with GitHub Actions" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.248.14"
$ws.Range("E2").Value = "  +1.68%  "

$ws.Range("D3").Value = "3.227.52"
$ws.Range("E3").Value = "  +1.24%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "606.64"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.67%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "155.75"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.78%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.05%  "

$ws.Range("D8").Value = "3.229.92"
$ws.Range("E8").Value = "  +1.37%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.536"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.72%  "

$ws.Range("E10").Value = "  -1.03%  "

$ws.Range("E11").Value = "  -0.38%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.512"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.75%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000272"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.09%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "38.98"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.95%  "

$ws.Range("D15").Value = "3.748.73"
$ws.Range("E15").Value = "  +1.00%  "

$ws.Range("B16").Value = "Polkadot"
$ws.Range("C16").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.49"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.91%  "

$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "66.362.33"
$ws.Range("E17").Value = "  +1.72%  "

$ws.Range("D18").Value = "3.236.11"
$ws.Range("E18").Value = "  +1.38%  "

$ws.Range("E19").Value = "  +0.63%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "514.44"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.24%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.83"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.91%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.741"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.50%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "15.37"
$ws.Range("D23").Style = "Normal"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.02"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.03%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.77"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.37%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.12%  "

$ws.Range("E27").Value = "  +4.18%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.31"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.00%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.25"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.65%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.90"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.15%  "

$ws.Range("E31").Value = "  +9.82%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.37"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.67%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.22"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.90%  "

$ws.Range("E34").Value = "  +0.11%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.70"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.46%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "55.67"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.17%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0926"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.61%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "493.48"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.66%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0424"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.21%  "

$ws.Range("E40").Value = "  -3.64%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.90"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.37%  "

$ws.Range("D42").Value = "3.037.68"
$ws.Range("E42").Value = "  -1.31%  "

$ws.Range("B43").Value = "TheGraph"
$ws.Range("C43").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.296"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.42%  "

$ws.Range("B44").Value = "Kaspa"
$ws.Range("C44").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.120"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.10%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.51"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.56%  "

$ws.Range("D46").Value = "0.0X0655"
$ws.Range("D46").Characters(4, 1).Text = [string][char]0x2083
$ws.Range("E46").Value = "  +7.13%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "29.33"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.37%  "

$ws.Range("E48").Value = "  +0.09%  "

$ws.Range("E49").Value = "  +0.10%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.34"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.70%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "120.13"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.95%  "
